$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the last four data rows (24-27): these used to belong to the
# previous run/block; they are now folded into the current run with new
# timestamps, phase names and durations. Row 27 is the "not storing image"
# statement (ResetTime) added at the tail of this run. ---

# Row 24
$ws.Cells.Item(24, 2).Value = 44270.476821837336
$ws.Cells.Item(24, 3).Value = "ShootTime"
$ws.Cells.Item(24, 4).Value = 2.0442542000000001

# Row 25
$ws.Cells.Item(25, 2).Value = 44270.476845528465
$ws.Cells.Item(25, 3).Value = "StoreTime"
$ws.Cells.Item(25, 4).Value = 0.63934100000000005

# Row 26
$ws.Cells.Item(26, 2).Value = 44270.47685413561
$ws.Cells.Item(26, 3).Value = "ConnTime"
$ws.Cells.Item(26, 4).Value = 11.1906836

# Row 27
$ws.Cells.Item(27, 2).Value = 44270.476854131281
$ws.Cells.Item(27, 3).Value = "ResetTime"
$ws.Cells.Item(27, 4).Value = 16.2050576

# --- Re-apply the DateTime number format across the whole run (B2:B27) so
# the block shares a single, freshly-minted style (mirrors the new
# cellXfs/borders entry the workbook gained for this run). ---
$ws.Range("B2:B27").NumberFormat = "m/d/yy h:mm"
$ws.Range("B2:B27").Borders.LineStyle = 1
